$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "Checklist" to "Session"
$ws.Name = "Session"

# Remove the last two data rows (old rows 4 and 5), shifting remaining rows up
$ws.Rows(4).Delete()
$ws.Rows(4).Delete()

# Update row 2 values (keep Student ID stored as text, like the rest of the sheet)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "231995"
$ws.Range("D2").Value = "09:50:06"
$ws.Range("E2").Value = "Scan"
$ws.Range("F2").Value = "5edfa2692bdacc5e6ee805c626c50cb44cebb065f092d9a1067d89f74dacd326"

# Update row 3 values
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "231997"
$ws.Range("E3").Value = "Scan"
$ws.Range("F3").Value = "5edfa2692bdacc5e6ee805c626c50cb44cebb065f092d9a1067d89f74dacd326"
